$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the sheet tab name (drop the trailing " 2023")
$ws.Name = "g3.5a"

# Add the new "Ano" column: copy the header formatting from the
# existing "Ordem" header cell (C1) so D1 picks up the same style,
# then set its text and the 2023 value for every data row.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Ano"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 4).Value = 2023
}
